# "upgrade left table until javakheti" - extend the Kazbegi remuneration
# table with the 2023 column (K), matching the formatting already used by
# the preceding year columns (B-J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the formatting (number format, borders, alignment, etc.) of the
# last existing year column (J) onto the new column (K).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New 2023 data.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1419.3
$ws.Range("K5").Value = 1069.8
$ws.Range("K6").Value = 1724.8
